$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sayfa1")

# Update D_min (C6) from 0.25 to 0.3
$ws.Range("C6").Value = 0.3

# Update L (C10) formula from 47*10^-6 to 160*10^-6
$ws.Range("C10").Formula = "=160*10^-6"

# Update the active selection to I10 to match the saved workbook state
$ws.Range("I10").Select()
